$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 70,2
$arr[0,0] = 3
$arr[0,1] = 1036
$arr[1,0] = 4
$arr[1,1] = 551
$arr[2,0] = 5
$arr[2,1] = 287
$arr[3,0] = 6
$arr[3,1] = 192
$arr[4,0] = 7
$arr[4,1] = 157
$arr[5,0] = 8
$arr[5,1] = 111
$arr[6,0] = 9
$arr[6,1] = 55
$arr[7,0] = 10
$arr[7,1] = 48
$arr[8,0] = 13
$arr[8,1] = 47
$arr[9,0] = 14
$arr[9,1] = 47
$arr[10,0] = 12
$arr[10,1] = 43
$arr[11,0] = 11
$arr[11,1] = 40
$arr[12,0] = 15
$arr[12,1] = 37
$arr[13,0] = 16
$arr[13,1] = 32
$arr[14,0] = 17
$arr[14,1] = 24
$arr[15,0] = 19
$arr[15,1] = 24
$arr[16,0] = 20
$arr[16,1] = 22
$arr[17,0] = 18
$arr[17,1] = 20
$arr[18,0] = 23
$arr[18,1] = 14
$arr[19,0] = 48
$arr[19,1] = 12
$arr[20,0] = 21
$arr[20,1] = 10
$arr[21,0] = 25
$arr[21,1] = 10
$arr[22,0] = 60
$arr[22,1] = 8
$arr[23,0] = 30
$arr[23,1] = 6
$arr[24,0] = 34
$arr[24,1] = 6
$arr[25,0] = 114
$arr[25,1] = 6
$arr[26,0] = 62
$arr[26,1] = 5
$arr[27,0] = 61
$arr[27,1] = 4
$arr[28,0] = 24
$arr[28,1] = 3
$arr[29,0] = 27
$arr[29,1] = 3
$arr[30,0] = 29
$arr[30,1] = 3
$arr[31,0] = 38
$arr[31,1] = 3
$arr[32,0] = 43
$arr[32,1] = 3
$arr[33,0] = 26
$arr[33,1] = 2
$arr[34,0] = 31
$arr[34,1] = 2
$arr[35,0] = 36
$arr[35,1] = 2
$arr[36,0] = 44
$arr[36,1] = 2
$arr[37,0] = 46
$arr[37,1] = 2
$arr[38,0] = 50
$arr[38,1] = 2
$arr[39,0] = 51
$arr[39,1] = 2
$arr[40,0] = 63
$arr[40,1] = 2
$arr[41,0] = 70
$arr[41,1] = 2
$arr[42,0] = 117
$arr[42,1] = 2
$arr[43,0] = 120
$arr[43,1] = 2
$arr[44,0] = 140
$arr[44,1] = 1
$arr[45,0] = 143
$arr[45,1] = 1
$arr[46,0] = 22
$arr[46,1] = 1
$arr[47,0] = 32
$arr[47,1] = 1
$arr[48,0] = 33
$arr[48,1] = 1
$arr[49,0] = 35
$arr[49,1] = 1
$arr[50,0] = 39
$arr[50,1] = 1
$arr[51,0] = 41
$arr[51,1] = 1
$arr[52,0] = 42
$arr[52,1] = 1
$arr[53,0] = 47
$arr[53,1] = 1
$arr[54,0] = 49
$arr[54,1] = 1
$arr[55,0] = 52
$arr[55,1] = 1
$arr[56,0] = 55
$arr[56,1] = 1
$arr[57,0] = 57
$arr[57,1] = 1
$arr[58,0] = 58
$arr[58,1] = 1
$arr[59,0] = 67
$arr[59,1] = 1
$arr[60,0] = 71
$arr[60,1] = 1
$arr[61,0] = 96
$arr[61,1] = 1
$arr[62,0] = 100
$arr[62,1] = 1
$arr[63,0] = 107
$arr[63,1] = 1
$arr[64,0] = 118
$arr[64,1] = 1
$arr[65,0] = 119
$arr[65,1] = 1
$arr[66,0] = 121
$arr[66,1] = 1
$arr[67,0] = 123
$arr[67,1] = 1
$arr[68,0] = 126
$arr[68,1] = 1
$arr[69,0] = 127
$arr[69,1] = 1

$ws.Range("A2:B71").Value = $arr

$ws.Columns.Item(2).ColumnWidth = 8.0

$ws.Sort.SortFields.Clear()
$f = $ws.Sort.SortFields.Add($ws.Range("B2:B71"), 0, 2)
$ws.Sort.SetRange($ws.Range("A2:B71"))
$ws.Sort.Apply()